# Caso de uso "Cobrar boleto en Tren" - split into cobroMolinete / devolucionMolinete,
# added Maquina.carga row, and updated BoletoTren.calcularValor description
# to mention the new integration points.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 7: MaquinaTren.cobrar -> MaquinaTren.cobroMolinete, now takes a tarjeta param ---
$ws.Range("C7").Value = "cobroMolinete"
$ws.Range("E7").Value = "tarjeta: Tarjeta"

# --- Row 8: MaquinaTren.devolverSaldo -> MaquinaTren.devolucionMolinete ---
$ws.Range("C8").Value = "devolucionMolinete"
$ws.Range("D8").Value = "devuelve la diferencia al pasar por el molinete de salida"
$ws.Range("F8").Value = "void"

# --- Row 9 gains the Maquina.carga entry that used to live on row 10, and
#     row 10 becomes blank again - move the content+format up one row ---
$ws.Range("B10:F10").Copy()
$ws.Range("B9:F9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteAll)
$ws.Range("F9").Value = "void"

$ws.Range("B11:F11").Copy()
$ws.Range("B10:F10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B10:F10").Value = ""
$ws.Application.CutCopyMode = $false

# --- Row 18: BoletoTren.calcularValor description now mentions the new methods, and
#     the description cell becomes merged D18:E18 like the rows above it ---
$ws.Range("D18").Value = "calcula el valor del boleto según la seccion – Integrado en cobroMolinete y devolucionMolinete (Tarjeta)"
$ws.Range("E18").Value = ""
$ws.Range("D18").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignLeft
$ws.Range("D18").VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter
$ws.Range("D18:E18").Merge()

# --- Selection moves to D18 as the last-active cell ---
$ws.Range("D18").Select()

$wb.Save()
